$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2014925373134328
$ws.Range("C2").Value = 0.5410447761194029
$ws.Range("J2").Value = 0.007462686567164179
$ws.Range("P2").Value = 0.1305970149253731
$ws.Range("S2").Value = 0.1194029850746269
$ws.Range("B3").Value = 0.02666666666666667
$ws.Range("C3").Value = 0.02
$ws.Range("J3").Value = 0.01333333333333333
$ws.Range("P3").Value = 0.7466666666666667
$ws.Range("S3").Value = 0.1933333333333333
$ws.Range("J4").Value = 0.08823529411764706
$ws.Range("P4").Value = 0.6764705882352942
$ws.Range("S4").Value = 0.2352941176470588
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.115
$ws.Range("D6").Value = 0.01
$ws.Range("F6").Value = 0.055
$ws.Range("J6").Value = 0.195
$ws.Range("O6").Value = 0.02
$ws.Range("Q6").Value = 0.12
$ws.Range("R6").Value = 0.08
$ws.Range("S6").Value = 0.405
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.0202020202020202
$ws.Range("F7").Value = 0.05555555555555555
$ws.Range("J7").Value = 0.1363636363636364
$ws.Range("O7").Value = 0.01515151515151515
$ws.Range("Q7").Value = 0.2121212121212121
$ws.Range("R7").Value = 0.0707070707070707
$ws.Range("S7").Value = 0.3787878787878788
$ws.Range("B8").Value = 0.07142857142857142
$ws.Range("D8").Value = 0.02232142857142857
$ws.Range("F8").Value = 0.06919642857142858
$ws.Range("J8").Value = 0.1316964285714286
$ws.Range("O8").Value = 0.01785714285714286
$ws.Range("Q8").Value = 0.1785714285714286
$ws.Range("R8").Value = 0.07142857142857142
$ws.Range("S8").Value = 0.4375
$ws.Range("B9").Value = 0.1055555555555556
$ws.Range("D9").Value = 0.01666666666666667
$ws.Range("F9").Value = 0.06111111111111111
$ws.Range("J9").Value = 0.1611111111111111
$ws.Range("O9").Value = 0.01666666666666667
$ws.Range("Q9").Value = 0.1833333333333333
$ws.Range("R9").Value = 0.05
$ws.Range("S9").Value = 0.4055555555555556
$ws.Range("B10").Value = 0.1022222222222222
$ws.Range("D10").Value = 0.01511111111111111
$ws.Range("E10").Value = 0.0008888888888888889
$ws.Range("F10").Value = 0.07377777777777778
$ws.Range("J10").Value = 0.1253333333333333
$ws.Range("O10").Value = 0.01777777777777778
$ws.Range("Q10").Value = 0.2302222222222222
$ws.Range("R10").Value = 0.06666666666666667
$ws.Range("S10").Value = 0.368
$ws.Range("F11").Value = 0.003401360544217687
$ws.Range("G11").Value = 0.173469387755102
$ws.Range("J11").Value = 0.07823129251700681
$ws.Range("K11").Value = 0.2346938775510204
$ws.Range("L11").Value = 0.5034013605442177
$ws.Range("S11").Value = 0.006802721088435374
$ws.Range("G12").Value = 0.7531645569620253
$ws.Range("J12").Value = 0.1835443037974684
$ws.Range("L12").Value = 0.05063291139240506
$ws.Range("S12").Value = 0.01265822784810127
$ws.Range("G13").Value = 0.8333333333333334
$ws.Range("J13").Value = 0.09523809523809523
$ws.Range("S13").Value = 0.07142857142857142
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.01904761904761905
$ws.Range("H15").Value = 0.1571428571428571
$ws.Range("I15").Value = 0.04761904761904762
$ws.Range("J15").Value = 0.3238095238095238
$ws.Range("K15").Value = 0.05714285714285714
$ws.Range("M15").Value = 0.009523809523809525
$ws.Range("O15").Value = 0.0761904761904762
$ws.Range("S15").Value = 0.3095238095238095
$ws.Range("F16").Value = 0.01818181818181818
$ws.Range("H16").Value = 0.1878787878787879
$ws.Range("I16").Value = 0.07878787878787878
$ws.Range("J16").Value = 0.4060606060606061
$ws.Range("K16").Value = 0.08484848484848485
$ws.Range("M16").Value = 0.01818181818181818
$ws.Range("O16").Value = 0.04242424242424243
$ws.Range("S16").Value = 0.1636363636363636
$ws.Range("F17").Value = 0.01847575057736721
$ws.Range("H17").Value = 0.1732101616628176
$ws.Range("I17").Value = 0.08545034642032333
$ws.Range("J17").Value = 0.4018475750577367
$ws.Range("K17").Value = 0.09699769053117784
$ws.Range("M17").Value = 0.0207852193995381
$ws.Range("O17").Value = 0.04849884526558892
$ws.Range("S17").Value = 0.1547344110854504
$ws.Range("F18").Value = 0.01351351351351351
$ws.Range("H18").Value = 0.2297297297297297
$ws.Range("I18").Value = 0.08108108108108109
$ws.Range("J18").Value = 0.3716216216216216
$ws.Range("K18").Value = 0.1013513513513514
$ws.Range("M18").Value = 0.006756756756756757
$ws.Range("O18").Value = 0.03378378378378379
$ws.Range("S18").Value = 0.1621621621621622
$ws.Range("F19").Value = 0.007818608287724784
$ws.Range("H19").Value = 0.2150117279124316
$ws.Range("I19").Value = 0.08444096950742767
$ws.Range("J19").Value = 0.3229085222830336
$ws.Range("K19").Value = 0.1110242376856919
$ws.Range("M19").Value = 0.0218921032056294
$ws.Range("N19").Value = 0.0007818608287724785
$ws.Range("O19").Value = 0.07818608287724785
$ws.Range("S19").Value = 0.1579358874120406
